$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1448.6666
$ws.Range("I11").Value = 1448.6666
$ws.Range("K11").Value = 1448.6666
$ws.Range("M11").Value = -1308.6666

$ws.Range("H17").Value = 2074.25
$ws.Range("J17").Value = 2074.25
$ws.Range("L17").Value = 6222.75
$ws.Range("N17").Value = -6558.75

$ws.Range("H40").Value = 3793.7568
$ws.Range("J40").Value = 4347.8076
$ws.Range("L40").Value = 4347.8076
$ws.Range("N40").Value = -4697.8076

$ws.Range("H100").Value = 1652.8
$ws.Range("I100").Value = 1507.4546
$ws.Range("J100").Value = 2052.5
$ws.Range("K100").Value = 1507.4546
$ws.Range("L100").Value = 2052.5
$ws.Range("M100").Value = -966.4546
$ws.Range("N100").Value = -3134.5

$ws.Range("H135").Value = 1831.0769
$ws.Range("I135").Value = 1310.5
$ws.Range("K135").Value = 11794.5
$ws.Range("M135").Value = -9259.5

$ws.Range("H137").Value = 2068.1724
$ws.Range("I137").Value = 2076.077
$ws.Range("K137").Value = 6228.231000000001
$ws.Range("M137").Value = -3678.231000000001

$ws.Range("H138").Value = 2729.8604
$ws.Range("J138").Value = 3001.4849
$ws.Range("L138").Value = 9004.4547
$ws.Range("N138").Value = -19284.4547

$ws.Range("H141").Value = 9329.666999999999
$ws.Range("I141").Value = 9500
$ws.Range("J141").Value = 9244.5
$ws.Range("K141").Value = 28500
$ws.Range("L141").Value = 27733.5
$ws.Range("M141").Value = -23320
$ws.Range("N141").Value = -38093.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10916.333
$ws.Range("I61").Value = 10916.333
$ws.Range("K61").Value = 10916.333
$ws.Range("M61").Value = -10704.333

$ws.Range("H74").Value = 2876.3704
$ws.Range("I74").Value = 1906.48
$ws.Range("K74").Value = 1906.48
$ws.Range("M74").Value = -1032.48

$ws.Range("H77").Value = 2876.3704
$ws.Range("I77").Value = 1906.48
$ws.Range("K77").Value = 9532.4
$ws.Range("M77").Value = -5164.4

$ws.Range("H136").Value = 10916.333
$ws.Range("I136").Value = 10916.333
$ws.Range("K136").Value = 32748.999
$ws.Range("M136").Value = -30198.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2647
$ws.Range("I94").Value = 2422
$ws.Range("K94").Value = 2422
$ws.Range("M94").Value = -1971

$ws.Range("H134").Value = 2187.3447
$ws.Range("I134").Value = 1903.1666
$ws.Range("K134").Value = 5709.4998
$ws.Range("M134").Value = -3174.4998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2295.8696
$ws.Range("I16").Value = 1364.0667
$ws.Range("J16").Value = 4043
$ws.Range("K16").Value = 1364.0667
$ws.Range("L16").Value = 4043
$ws.Range("M16").Value = -1077.0667
$ws.Range("N16").Value = -4617

$ws.Range("H53").Value = 29437.25
$ws.Range("J53").Value = 29437.25
$ws.Range("L53").Value = 29437.25
$ws.Range("N53").Value = -30651.25

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H105").Value = 644.8570999999999
$ws.Range("I105").Value = 643.4
$ws.Range("K105").Value = 643.4
$ws.Range("M105").Value = 1103.6

$ws.Range("H113").Value = 2295.8696
$ws.Range("I113").Value = 1364.0667
$ws.Range("J113").Value = 4043
$ws.Range("K113").Value = 1364.0667
$ws.Range("L113").Value = 4043
$ws.Range("M113").Value = 805.9332999999999
$ws.Range("N113").Value = -8383

$ws.Range("H132").Value = 3179.25
$ws.Range("I132").Value = 2081.0667
$ws.Range("J132").Value = 6473.8
$ws.Range("K132").Value = 6243.2001
$ws.Range("L132").Value = 19421.4
$ws.Range("M132").Value = -3713.2001
$ws.Range("N132").Value = -24481.4

$ws.Range("H134").Value = 2527
$ws.Range("I134").Value = 1942.7333
$ws.Range("J134").Value = 4279.8
$ws.Range("K134").Value = 5828.199900000001
$ws.Range("L134").Value = 12839.4
$ws.Range("M134").Value = -3293.199900000001
$ws.Range("N134").Value = -17909.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 9228.143
$ws.Range("I122").Value = 5699.5
$ws.Range("K122").Value = 51295.5
$ws.Range("M122").Value = -48845.5

$ws.Range("H132").Value = 38463770
$ws.Range("I132").Value = 166667820
$ws.Range("K132").Value = 1500010380
$ws.Range("M132").Value = -1500007850

$ws.Range("H137").Value = 7471.125
$ws.Range("I137").Value = 5653.8
$ws.Range("J137").Value = 10500
$ws.Range("K137").Value = 16961.4
$ws.Range("L137").Value = 31500
$ws.Range("M137").Value = -11861.4
$ws.Range("N137").Value = -41700

$ws.Range("H140").Value = 1630.4286
$ws.Range("I140").Value = 1630.4286
$ws.Range("K140").Value = 4891.2858
$ws.Range("M140").Value = 288.7142000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3315.8125
$ws.Range("I132").Value = 2280.1538
$ws.Range("K132").Value = 6840.4614
$ws.Range("M132").Value = -4310.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6613.8125
$ws.Range("I100").Value = 5711.222
$ws.Range("J100").Value = 7774.2856
$ws.Range("K100").Value = 5711.222
$ws.Range("L100").Value = 7774.2856
$ws.Range("M100").Value = -5170.222
$ws.Range("N100").Value = -8856.285599999999

$ws.Range("H134").Value = 97851.336
$ws.Range("J134").Value = 97851.336
$ws.Range("L134").Value = 97851.336
$ws.Range("N134").Value = -107991.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6661.6665
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 6661.6665
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 6661.6665
$ws.Range("N62").Value = -7909.6665
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 6661.6665
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 6661.6665
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 33308.3325
$ws.Range("N65").Value = -39548.3325
$ws.Range("M65").ClearContents()

$ws.Range("H70").Value = 40105
$ws.Range("J70").Value = 40105
$ws.Range("L70").Value = 40105
$ws.Range("N70").Value = -40735

$ws.Range("H73").Value = 40105
$ws.Range("J73").Value = 40105
$ws.Range("L73").Value = 40105
$ws.Range("N73").Value = -42289
